$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3 all share this string)
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $targets = @()
    foreach ($cell in $used.Cells) {
        # NOTE: cast to [string] explicitly -- this host's Range.Text can come
        # back as a typed Boolean for cells like "True"/"False", and a bare
        # `-eq "Ready for handoff"` comparison against a boolean coerces the
        # string operand to bool (truthy) instead of comparing text, which
        # would wrongly "match" unrelated True/False cells.
        $text = [string]$cell.Text
        if ($text -eq "Ready for handoff") {
            $targets += $cell.Address()
        }
    }
    foreach ($addr in $targets) {
        $ws.Range($addr).Value = "In Translation"
    }
}

# ---------------------------------------------------------------------------
# 2) Narrower "Status" columns now that the label text is shorter.
#    Target stored widths (~13.41 chars) sit between this host's pixel-grid
#    steps, so the closest representable ColumnWidth is used.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
